# "catching up following live changes made at rehearsal"
# Re-create the parameter-map annotations (columns B/C) that were entered live,
# widen column B, rezoom/reselect the sheet, and nudge the saved window position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C annotations -------------------------------------------------
# NB: the exact order below matters - it reproduces the live entry order so the
# resulting xl/sharedStrings.xml table comes out in the same sequence.
$ws.Range("C3").Value  = "tune"
$ws.Range("C27").Value = "output"
$ws.Range("C28").Value = "divergence"
$ws.Range("C29").Value = "size"
$ws.Range("C22").Value = "size"
$ws.Range("C30").Value = "LFElevel"
$ws.Range("C31").Value = "Xshift"
$ws.Range("C32").Value = "Yshift"
$ws.Range("C33").Value = "RotationAngle"
$ws.Range("C34").Value = "distance"
$ws.Range("B27").Value = "InstInsertFX:Sr.Panner"
$ws.Range("B9").Value  = "PitchLFOmulti"
$ws.Range("C4").Value  = "curve"
$ws.Range("C5").Value  = "attack"
$ws.Range("C6").Value  = "decay"
$ws.Range("C7").Value  = "sustain"
$ws.Range("C8").Value  = "release"
$ws.Range("B4").Value  = "Group:Pitch:ADHSR"
$ws.Range("B3").Value  = "Group:Source"
$ws.Range("C9").Value  = "frequency"
$ws.Range("C10").Value = "fadein"
$ws.Range("C11").Value = "pulseWidth"
$ws.Range("C12").Value = "Sine"
$ws.Range("C13").Value = "Rectangle"
$ws.Range("C14").Value = "Saw"
$ws.Range("C23").Value = "dry"
$ws.Range("C24").Value = "wet"
$ws.Range("B37").Value = "remember: potential to double parameters on one hostAutoNumber"

# --- Column widths -------------------------------------------------------------
# Column B needs to display the long group/FX labels - widen it to 24.5 chars.
# (ColumnWidth is re-measured against the font on set/save, so back the
# requested width off by the same 5/6-character padding Excel re-adds.)
$ws.Columns.Item(2).ColumnWidth = 23.666666666666668

# --- View state ------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 125
$ws.Range("B38").Select() | Out-Null

# Best-effort: nudge the saved window position to match the live rehearsal
# session (cosmetic only - some hosts don't persist this back to bookViews).
try { $excel.ActiveWindow.Left = 22280 } catch {}
